# Add a right-aligned URL textbox (drawn as a Rectangle autoshape, matching
# the authored shape's name/geometry) to the title slide (slide 1).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Target geometry is expressed in EMU in the OOXML; PowerPoint's Shapes.AddShape
# takes Left/Top/Width/Height in points, so convert (1 pt = 12700 EMU).
$emuPerPoint = 12700
$left   = 1507067   / $emuPerPoint
$top    = 5980102   / $emuPerPoint
$width  = 10410506  / $emuPerPoint
$height = 646331    / $emuPerPoint

# msoShapeRectangle = 1
$shp = $s.Shapes.AddShape(1, $left, $top, $width, $height)

# Wrap text within the shape and let the shape grow/shrink to fit the text
# (<a:bodyPr wrap="square"><a:spAutoFit/></a:bodyPr>).
# msoTrue = -1
$shp.TextFrame.WordWrap = -1
# ppAutoSizeShapeToFitText = 1
$shp.TextFrame.AutoSize = 1

$tr = $shp.TextFrame.TextRange
$tr.Text = "https://github.com/ericburcham/anti-patterns"
$tr.Font.Size = 36
$tr.Font.Bold = $true
# ppAlignRight = 3
$tr.ParagraphFormat.Alignment = 3
